$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v0 = @'
Microsoft Azure
'@
$ws.Range("A3").Value = $v0

$v1 = @'
This is a cloud-based platform primarily used for SSL authentication and access management. It supports user licensing, email services, collaboration tools, and document management functionalities.
'@
$ws.Range("B3").Value = $v1

$v2 = @'
The client utilizes this system for SSL authentication and access management across different systems and for core Microsoft 365 services, including user licensing, collaboration through SharePoint, email via Exchange, and communication using Teams.
'@
$ws.Range("C3").Value = $v2

$v3 = @'
The system is currently administered by Ralph Vaccaro, with oversight expected to transition to Nicole Tai (Vice President) upon her joining.
'@
$ws.Range("D3").Value = $v3

$v4 = @'
Access provisioning is initiated by HR sending a message with the new hire's details, including name, department, position, and start date, to the IT team, who manually create accounts prior to the start date. Elevated access requests require submission by the hiring manager or area lead and approval from IT leadership or another department head before access is granted.
'@
$ws.Range("E3").Value = $v4

$v5 = @'
Access is removed upon termination or role change when HR or the hiring manager sends an email specifying the change and effective date. The IT team deactivates access starting with Azure, which cascades to SSO-integrated systems, while non-SSO systems are handled independently.
'@
$ws.Range("F3").Value = $v5

$v6 = @'
Access is configured using a role-based model, where predefined roles such as regular user and elevated access (admin) are used, with elevated access granting full administrative privileges; currently, only one individual and a dormant backup account hold elevated access.
'@
$ws.Range("G3").Value = $v6

$v7 = @'
No; Only vendor-provided roles exist, and the client cannot alter role definitions without vendor assistance.
'@
$ws.Range("H3").Value = $v7

$v8 = @'
No; management does not perform periodic reviews of roles and permissions.
'@
$ws.Range("I3").Value = $v8

$v9 = @'
Yes; Users with privileged access include Ralph Vaccaro with active admin access and a dormant IT Admin account used as a backup.
'@
$ws.Range("J3").Value = $v9

$v10 = @'
Yes; There is an interactive IT Admin account used as a backup, and additional generic accounts exist for specific functionalities, but they do not have elevated access.
'@
$ws.Range("K3").Value = $v10

$v11 = @'
The credentials for the non-human IT admin account are memorized by Ralph Vaccaro, with no formal documentation or storage method.
'@
$ws.Range("L3").Value = $v11

$v12 = @'
The password for the IT Admin account is known only to Ralph Vaccaro and stored in his memory, while other non-human accounts with standard access are configured as regular users without elevated access.
'@
$ws.Range("M3").Value = $v12

$v13 = @'
No; Management does not currently perform periodic access reviews. Elevated access in Azure has not been reviewed as of today.
'@
$ws.Range("N3").Value = $v13

$v14 = @'
No, the system does not have activity logging capabilities or audit trail functionality.
'@
$ws.Range("O3").Value = $v14

$v15 = @'
No; periodic reviews of user activity, roles, permissions, or elevated access are not currently performed, and activity logging or tracking functionality is not actively utilized.
'@
$ws.Range("P3").Value = $v15

$v16 = @'
The system uses Azure SSO with MFA for all users, and privileged accounts include a backup non-human IT Admin account for emergency access.
'@
$ws.Range("Q3").Value = $v16

$v17 = @'
No; The client does not perform periodic reviews of the systems authentication configurations.
'@
$ws.Range("R3").Value = $v17

$v18 = @'
Management can modify permissions on certain restricted areas within Azure. No other types of changes, such as configuration, workflow, or code changes, are performed by management.
'@
$ws.Range("S3").Value = $v18

$v19 = @'
Access to make changes is restricted to Ralph Vaccaro, who holds an admin role with elevated access, and the IT Admin Account, a dormant non-human account used as a backup. Ralph is the only individual with knowledge of the IT Admin Account password, ensuring controlled access to modify the system.
'@
$ws.Range("T3").Value = $v19

$v20 = @'
Management does not have any separate environments for this system.
'@
$ws.Range("U3").Value = $v20

$v21 = @'
There is no formal change management process in place. Changes are requested verbally or via email, and the IT team makes changes directly in the production environment without formal approvals or testing. There is no documentation of changes, and no review process exists to validate changes were appropriate.
'@
$ws.Range("V3").Value = $v21

$v22 = @'
No, the vendor manages updates, patches, and bug fixes for the SaaS system. The organization does not perform independent testing or validation of these updates and relies entirely on the vendor for deployment and oversight.
'@
$ws.Range("W3").Value = $v22

$v23 = @'
N/A - This information was not discussed in the walkthrough meeting transcript.
'@
$ws.Range("X3").Value = $v23

$v24 = @'
No periodic review of changes is performed.
'@
$ws.Range("Y3").Value = $v24

$v25 = @'
No automated jobs or interfaces are currently implemented for this system.
'@
$ws.Range("Z3").Value = $v25

$v26 = @'
The automated jobs are configured using native system functionality and are not scheduled to run, as there are no automated jobs or interfaces currently in use between Azure and other systems.
'@
$ws.Range("AA3").Value = $v26

$v27 = @'
N/A - This information was not discussed in the walkthrough meeting transcript.
'@
$ws.Range("AB3").Value = $v27

$v28 = @'
The system utilizes Microsoft Azure for data storage, which is vendor-managed.
'@
$ws.Range("AC3").Value = $v28

$v29 = @'
Backups are not currently performed independently, but implementation of Druva for 365 backups is in progress.
'@
$ws.Range("AD3").Value = $v29

$v30 = @'
N/A - All backup types and strategies are managed by the vendor as part of their SaaS service.
'@
$ws.Range("AE3").Value = $v30

$v31 = @'
N/A - the vendor is responsible for monitoring and resolving any backup failures.
'@
$ws.Range("AF3").Value = $v31

$v32 = @'
No, management does not perform regular SOC report reviews.
'@
$ws.Range("AG3").Value = $v32
